# AppSG_cea_09. 3-tier configuration.xlsx — "Add files via upload" edit
#
# The HTTPS/HTTP row (row 4, "0.0.0.0/0" target / Outbound) had its TCP Port
# cell typed as a bare number (80443). The re-uploaded version fixes this to
# the intended comma-separated text "80, 443" (matching the sibling cells
# E3/E5/E6 convention of listing ports), and the author's cursor/selection
# ended up parked on E5 when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SecurityGroupRule_List (2)")

# E4 was a numeric literal (80443); replace it with the text "80, 443".
# Assigning a string through .Value (without touching NumberFormat) makes
# Excel store it as an inline/shared text string while leaving the cell's
# existing style (s="17") untouched, exactly like a user typing text over
# a numeric-formatted cell.
$ws.Range("E4").Value = "80, 443"

# Match the saved selection/view state (cursor on E5, scrolled so row 2 is
# at the top of the pane) from the re-uploaded workbook.
$ws.Range("E5").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
